# project_tracker.xlsx update - "data management files added to dir"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Detector unit (rows 5-10) -------------------------------------------------
# Row 5 (Backend comms framework): now COMPLETE (was IN PROGRESS)
$ws.Range("H5").Value = "COMPLETE"

# Row 7 (Assemble RC components): remark updated
$ws.Range("F7").Value = "Using PQ-2 model (once built)"

# Row 9 (Test power pi from 5V pin on FC): remark updated
$ws.Range("F9").Value = "Perhaps better to use 18650 with shield? Compare FC pin with this"

# Row 9 network/interface column (Test issue commands...): now COMPLETE (was IN PROGRESS)
$ws.Range("H9").Value = "COMPLETE"

# --- Rover (rows 11-16) --------------------------------------------------------
# Row 11 (Telemetry over 4G): now IN PROGRESS, remark updated
$ws.Range("E11").Value = "IN PROGRESS"
$ws.Range("F11").Value = "Seems to work fine over VPN"

# --- Print / Assemble section (rows 21-22) -------------------------------------
# Row 21 (Design servo mount and pi/camera case): now IN PROGRESS, remark updated
$ws.Range("B21").Value = "IN PROGRESS"
$ws.Range("C21").Value = "Have designed a payload container for testing for now, later need to design with servo mounts & rotating camera"

# Row 22 (Print): now IN PROGRESS with a remark
$ws.Range("B22").Value = "IN PROGRESS"
$ws.Range("C22").Value = "Printing provisional/testing Rpi payload container"

# Row heights for the newly-populated / re-wrapped rows
$ws.Rows.Item(22).RowHeight = 46.25
$ws.Rows.Item(23).RowHeight = 35.8

# Restore the active selection to H9
[void]$ws.Range("H9").Select()
